$wb = $excel.ActiveWorkbook

$wsProfile = $wb.Worksheets.Item("CreateModifyDeleteProfile")
$wsCC      = $wb.Worksheets.Item("AddModifyDeleteCC")
$wsACH     = $wb.Worksheets.Item("AddModifyDeleteACH")

$wsProfile.Range("B2").Value = "Thu Jun 19 19:15:11 IST 2025"
$wsCC.Range("B2").Value      = "Thu Jun 19 18:06:58 IST 2025"
$wsACH.Range("B2").Value     = "Thu Jun 19 18:03:23 IST 2025"
$wsACH.Range("B3").Value     = "Thu Jun 19 18:04:38 IST 2025"
$wsACH.Range("B4").Value     = "Thu Jun 19 18:05:48 IST 2025"
